$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 8 new rows before the old "Total" row (row 53), pushing it (and the
#    notes below it) down to make room for a new journal-entry block, mirroring
#    the existing "Copier/Inserer un jour supplementaire" pattern used for every
#    other day block on this sheet.
$ws.Rows("53:60").Insert()

# 2. Seed the new block by copying the previous complete day-block (rows 45-52)
#    down into the newly inserted rows - this carries over the row heights,
#    merged cell B46:B52 shape, borders and styles in one shot.
$ws.Range("A45:E45").Copy($ws.Range("A53:E53"))
$ws.Range("A46:E46").Copy($ws.Range("A54:E54"))
$ws.Range("A47:E52").Copy($ws.Range("A55:E60"))

# 3. A couple of the freshly-typed rows mix the "blank row" label/notes style
#    with the "value" style for the duration column - reproduce that exactly.
$ws.Range("A50").Copy($ws.Range("A56"))
$ws.Range("A50").Copy($ws.Range("A57"))
$ws.Range("C47").Copy($ws.Range("C56"))
$ws.Range("C47").Copy($ws.Range("C57"))
$ws.Range("D50").Copy($ws.Range("D56"))
$ws.Range("D50").Copy($ws.Range("D57"))
$ws.Range("E47").Copy($ws.Range("E56"))
$ws.Range("E47").Copy($ws.Range("E57"))

# 4. Fill in the new journal entry for 2024-05-06 (serial 45422) - a new
#    "Affichage Sondage" sub-task of the "CRUD" task.
$ws.Range("A54").Value = "CRUD"
$ws.Range("B54").Value = 45422
$ws.Range("C54").Value = 300
$ws.Range("D54").Value = "Réaliser la possibilité aux utilisateurs de créer des sondages, ceux-ci gènère directement les questions et les réponses"
$ws.Range("E54").Value = "Cela m'a pris plus de temps que prévu"

$ws.Range("A55").Value = "Affichage Sondage"
$ws.Range("C55").Value = 50
$ws.Range("D55").Value = "Création de l'affichage des sondages (En cours)"

$ws.Range("A56").Value = "Pause"
$ws.Range("C56").Value = 5
$ws.Range("D56").Value = "Pause donnée par le surveillant"

$ws.Range("A57").Value = "JDT"
$ws.Range("C57").Value = 5
$ws.Range("D57").Value = "Remplire le journal de travail"

# 5. Extend the daily-total formula so it also covers the new block, now that
#    the total row has shifted from row 53 down to row 61.
$ws.Range("C61").Formula = "=MROUND(SUM(C7:C36,C54:C60) /60,0.2)"

# 6. Keep the print area in sync with the sheet's new extent.
$ws.PageSetup.PrintArea = "A1:E61"

$ws.Calculate()
